$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.53"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.90"
$ws.Range("D3").ClearFormats()
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.178"
$ws.Range("D4").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.474"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.149"
$ws.Range("D7").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8142"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8539"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1375"
$ws.Range("D10").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06953"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03204"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCXBestin24h"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02886"
$ws.Range("D13").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09383"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.826"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001537"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04704"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005988"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006158"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004780"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00008498"
$ws.Range("D22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.532"
$ws.Range("D23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.157"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3202"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1330"
$ws.Range("D26").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1358"
$ws.Range("D27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002330"
$ws.Range("D28").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03697"
$ws.Range("D40").ClearFormats()
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1055"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.002205"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003067"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "42KickTokenKICK"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007808"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005487"
$ws.Range("D45").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002556"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("D49").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("D50").ClearFormats()
